# edit.ps1
# Applies the "Finished Slide show 4 and fix user manual" edits:
#  1. Slide 2 ("Roles"): Muhammed Carrim box - turn the trailing ": " into
#     ":" and add two new (non-bold) bullet lines describing his work.
#  2. Slide 4 ("Impressive aspect"): fill in the previously empty content
#     placeholder with four bullet points.
#  3. Slide 6 ("List of what is left"): fill in the previously empty
#     content placeholder with four bullet points.
#  4. Slide 7 ("Stand-out extension"): fill in the previously empty content
#     placeholder with four bullet points plus a trailing, unbulleted blank
#     paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 2 - "Muhammed Carrim" rounded rectangle (Rounded Rectangle 5)
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$carrimShape = $slide2.Shapes.Item(4)

# Re-assert the existing heading text so it keeps its current run
# formatting, just drop the trailing space after the colon.
$carrimShape.TextFrame.TextRange.Text = "Muhammed Carrim:"
# Now grow the text to include the new bullet lines; because a run already
# exists with lang="en-US"/sz="1400"/b="1" the new paragraphs inherit it.
$carrimShape.TextFrame.TextRange.Text = "Muhammed Carrim:" + "`r" + "DDOS Attack Test," + "`r" + "Helped with Load Balancing"

# The two new lines are not bold in the target, so clear bold on them.
$carrimShape.TextFrame.TextRange.Paragraphs(2, 1).Font.Bold = $false
$carrimShape.TextFrame.TextRange.Paragraphs(3, 1).Font.Bold = $false

# ---------------------------------------------------------------------------
# 2. Slide 4 - "Impressive aspect" content placeholder
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$impressiveShape = $slide4.Shapes.Item(2)

$impressiveLine1 = "Can protect against DDOS attack instead of DOS only."
$impressiveRest = @(
    "Heatmap of Incoming Traffic.",
    "Different Load Balancing Algorithms",
    "Shell file to install all dependancies"
)

# Seed a real run (with lang="en-US") by writing the first line alone first,
# then grow to the full, multi-paragraph text so every new paragraph
# inherits that run's language/formatting.
$impressiveShape.TextFrame.TextRange.Text = $impressiveLine1
$impressiveShape.TextFrame.TextRange.Text = $impressiveLine1 + "`r" + ($impressiveRest -join "`r")

# ---------------------------------------------------------------------------
# 3. Slide 6 - "List of what is left" content placeholder
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$listLeftShape = $slide6.Shapes.Item(2)

$listLeftLine1 = "Optimizing database connection"
$listLeftRest = @(
    "Adding triggers to database",
    "Exponential counter for blacklisting",
    "Improving interface to backend compatibility"
)

$listLeftShape.TextFrame.TextRange.Text = $listLeftLine1
$listLeftShape.TextFrame.TextRange.Text = $listLeftLine1 + "`r" + ($listLeftRest -join "`r")

# ---------------------------------------------------------------------------
# 4. Slide 7 - "Stand-out extension" content placeholder
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$standoutShape = $slide7.Shapes.Item(2)

$standoutLine1 = "Logging provides capabilities for forensic investigation into possible attacks"
$standoutRest = @(
    "Can protect any backend server that is registered with the device",
    "Possibility to defend against SYN attacks from number of connections implementation",
    "Documentation is user friendly "
)

$standoutShape.TextFrame.TextRange.Text = $standoutLine1
# Trailing "`r" leaves one extra, empty paragraph at the end (matching the
# blank, bullet-less paragraph in the target).
$standoutShape.TextFrame.TextRange.Text = $standoutLine1 + "`r" + ($standoutRest -join "`r") + "`r"

# That final blank paragraph has no bullet in the target.
$lastParaIndex = $standoutShape.TextFrame.TextRange.Paragraphs().Count
$standoutShape.TextFrame.TextRange.Paragraphs($lastParaIndex, 1).ParagraphFormat.Bullet.Type = 0
